$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts rows 13-23 down to 14-24),
# preserving row heights and existing cell formatting automatically.
$ws.Rows("13:13").Insert()

# New row 13: Docentes responsaveis (B13/C13)
$docentes = @'
8554681 - Pedro Felipe Arce Castillo
'@
$ws.Range("B13").Value = $docentes
$ws.Range("C13").Value = $docentes

# Update Objetivos text (B10/C10)
$objetivos = @'
Ao final do curso os estudantes deverão: - Compreender os aspectos mássicos, energéticos e entrópicos, envolvendo sistemas termodinâmicos abertos e fechados; - Dominar e ser capaz de fazer predições básicas de propriedades termodinâmicas, usando equações cúbicas de estado e relações termodinâmicas; - Desenvolver uma metodologia para poder solucionar os problemas de engenharia, nos aspectos termodinâmicos; Dominar o uso de tabelas de propriedades termodinâmicas;
'@
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# Update Programa resumido text (B14/C14)
$programaResumido = @'
A primeira Lei da Termodinâmica. Efeitos Térmicos. A segunda lei da Termodinâmica. Propriedades termodinâmicas dos fluidos. Termodinâmica de processos de escoamento. Produção de potencia a partir de calor. Refrigeração e liquefação
'@
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido

# Update Programa (full syllabus) text (B16/C16)
$programa = @'
1  A primeira Lei da Termodinâmica
     1.1- Energia interna
     1.2- Estado termodinâmico e funções de estado
     1.3- Entalpia
     1.4- A regra das fases
     1.5- O processo reversível
     1.6- Processos a volume constante e a pressão constante
     1.7- Capacidade calorífica
2  Efeitos Térmicos 
     2.1- Calores Latentes de Substâncias Puras.
     2.2- Calor de Reação Padrão
     2.3- Calor Padrão de Formação
     2.4- Calor Padrão de Combustão
     2.5- O processo reversível
     2.6- A variação da entalpia com a Temperatura
3- A segunda lei da Termodinâmica
    3.1- Enunciados da lei
    3.2- Máquinas térmicas
    3.3- Escalas de temperaturas termodinâmicas
    3.4- Entropia 
    3.5- Variações da entropia de um gás ideal
    3.6- A terceira lei da termodinâmica
4- Produção de potencia a partir de calor
   4.1- A planta de potencia a vapor (maquina a vapor)
   4.2- Motores de combustão interna
   4.3- O motor Otto
   4.4- O motor Diesel
   4.5- A planta de potencia com turbina a gás
5- Refrigeração e liquefação
    5.1- O refrigerador de Carnot
    5.2- O ciclo com compresso a vapor
    5.3- Comparação de ciclos de refrigeração
    5.4- Refrigeração por absorção
    5.5- A bomba a calor
    6.6- Processos de liquefação
6- Termodinâmica de soluções
    6.1- Relações fundamentais entre propriedades
    6.2- O potencial químico 
    6.3- Fugacidade e coeficiente de fugacidade
    6.4- A solução Ideal
    6.5- Modelos para a energia de Gibbs
    6.6- Propriedades de mistura
    6.7- Efeitos térmicos em processos de mistura
7- Equilíbrio de fases
    7.1- Equilíbrio e estabilidade
    7.2- Equilíbrio líquido-líquido
    7.3- Equilíbrio líquido-líquido-vapor
    7.4- Equilíbrio sólido-líquido
    7.5- Equilíbrio sólido-vapor
    7.6- Equilíbrio na adsorção de gases em sólidos
'@
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# Update Metodo text (B19/C19)
$metodo = @'
2 provas escritas
'@
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# Update Criterio text (B20/C20)
$criterio = @'
serão avaliados os conteúdos discutidos em sala e constantes da ementa do curso. A média da disciplina será a média aritmética das duas provas.
'@
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# Update Norma de recuperacao text (B21/C21)
$norma = @'
prova escrita com conteúdo de todo o semestre
'@
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# Update Bibliografia text (B22/C22)
$bibliografia = @'
1)Smith, J.M.; Van Ness, H.C.; Abott, M. M.  Introdução à Termodinâmica da Engenharia Química. 7ª ed.  ISBN 978-85-216-1553-8, LTC editora, 2007.
2)Koretsky, M. D. Termodinâmica para Engenharia Química, 1ª ed.  ISBN 978-85-216-1530-9, LTC editora, 2007.
3)Terron, L. R. Termodinâmica Química Aplicada. 1ª ed.  ISBN 978-85-204-2082-9, Editora Manole Ltda, 2009.
4)Moran, M. J.; Shapiro, H. N. Princípios de Termodinâmica para Engenharia - 1ª ed.  ISBN 978-85-216-1689-4, LTC editora, 2009.
5)Van Wilen, J. Sonntag, Richard. E. Fundamentos da Termodinâmica Clássica  6ª Edição  2004
6)Sandler, S. I., Chemical and Engineering Thermodynamics, 3rd ed., John Wiley & Sons, 1999
'@
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia

